$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header style (bold, centered, bordered) from H1 onto the new I1:J1 headers
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# New header labels
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# New data values for columns I (I0) and J (IF), rows 2-40
$data = @(
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(6,7),
    @(6,6),
    @(8,9),
    @(9,9),
    @(1,2),
    @(1,4),
    @(1,6),
    @(1,6),
    @(1,2),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,7),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,6),
    @(1,6),
    @(1,5),
    @(1,2),
    @(1,6),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,5),
    @(1,4),
    @(3,5),
    @(4,6),
    @(3,5),
    @(1,2)
)

$r = 2
foreach ($pair in $data) {
    $ws.Cells.Item($r, 9).Value = $pair[0]
    $ws.Cells.Item($r, 10).Value = $pair[1]
    $r++
}
